# Update "想去人数" (want-to-go count) values in column F across sheets
# to reflect the latest scrape output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12887
$ws1.Range("F3").Value = 634
$ws1.Range("F6").Value = 329
$ws1.Range("F7").Value = 410
$ws1.Range("F9").Value = 12999
$ws1.Range("F10").Value = 45
$ws1.Range("F11").Value = 29
$ws1.Range("F12").Value = 5281
$ws1.Range("F13").Value = 550
$ws1.Range("F15").Value = 17
$ws1.Range("F16").Value = 34
$ws1.Range("F17").Value = 1202
$ws1.Range("F18").Value = 43
$ws1.Range("F19").Value = 136
$ws1.Range("F20").Value = 685
$ws1.Range("F21").Value = 2861
$ws1.Range("F22").Value = 6207
$ws1.Range("F23").Value = 1165
$ws1.Range("F24").Value = 3637
$ws1.Range("F26").Value = 48

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 28
$ws2.Range("F3").Value = 10

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12887
$ws4.Range("F3").Value = 634
$ws4.Range("F6").Value = 329
$ws4.Range("F7").Value = 28
$ws4.Range("F8").Value = 410
$ws4.Range("F10").Value = 12999
$ws4.Range("F11").Value = 45
$ws4.Range("F12").Value = 29
$ws4.Range("F13").Value = 5281
$ws4.Range("F14").Value = 550
$ws4.Range("F16").Value = 17
$ws4.Range("F17").Value = 34
$ws4.Range("F18").Value = 1202
$ws4.Range("F19").Value = 43
$ws4.Range("F20").Value = 136
$ws4.Range("F21").Value = 685
$ws4.Range("F22").Value = 2861
$ws4.Range("F23").Value = 10
$ws4.Range("F24").Value = 6207
$ws4.Range("F25").Value = 1165
$ws4.Range("F26").Value = 3637
$ws4.Range("F28").Value = 48
